$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right before row 938; this pushes the existing
# rows 938-1031 down to 941-1034 (preserving their values/styles).
$ws.Range("A938:A940").EntireRow.Insert()

# Constant columns shared by every data row in this sheet.
$mercadoId = 5
$mercado   = "Macroferia Regional de Talca"
$region    = "Maule"
$codreg    = 7
$catId     = 100112002
$categoria = "Pimiento"
$clasif    = "Hortaliza"

# Row 938
$ws.Range("A938").Value = $mercadoId
$ws.Range("B938").Value = $mercado
$ws.Range("C938").Value = $region
$ws.Range("D938").Value = 45194
$ws.Range("E938").Value = $codreg
$ws.Range("F938").Value = $catId
$ws.Range("G938").Value = $categoria
$ws.Range("H938").Value = "Zafiro rojo"
$ws.Range("I938").Value = "Primera"
$ws.Range("J938").Value = 100
$ws.Range("K938").Value = 45000
$ws.Range("L938").Value = 45000
$ws.Range("M938").Value = 45000
$ws.Range("N938").Value = "`$/caja 15 kilos"
$ws.Range("O938").Value = "Región de Arica y Parinacota"
$ws.Range("P938").Value = 3000
$ws.Range("Q938").Value = 15
$ws.Range("R938").Value = $clasif

# Row 939
$ws.Range("A939").Value = $mercadoId
$ws.Range("B939").Value = $mercado
$ws.Range("C939").Value = $region
$ws.Range("D939").Value = 45194
$ws.Range("E939").Value = $codreg
$ws.Range("F939").Value = $catId
$ws.Range("G939").Value = $categoria
$ws.Range("H939").Value = "Zafiro rojo"
$ws.Range("I939").Value = "Segunda"
$ws.Range("J939").Value = 100
$ws.Range("K939").Value = 40000
$ws.Range("L939").Value = 40000
$ws.Range("M939").Value = 40000
$ws.Range("N939").Value = "`$/caja 15 kilos"
$ws.Range("O939").Value = "Región de Arica y Parinacota"
$ws.Range("P939").Value = 2667
$ws.Range("Q939").Value = 15
$ws.Range("R939").Value = $clasif

# Row 940
$ws.Range("A940").Value = $mercadoId
$ws.Range("B940").Value = $mercado
$ws.Range("C940").Value = $region
$ws.Range("D940").Value = 45194
$ws.Range("E940").Value = $codreg
$ws.Range("F940").Value = $catId
$ws.Range("G940").Value = $categoria
$ws.Range("H940").Value = "Zafiro verde"
$ws.Range("I940").Value = "Primera"
$ws.Range("J940").Value = 100
$ws.Range("K940").Value = 35000
$ws.Range("L940").Value = 35000
$ws.Range("M940").Value = 35000
$ws.Range("N940").Value = "`$/caja 15 kilos"
$ws.Range("O940").Value = "Región de Arica y Parinacota"
$ws.Range("P940").Value = 2333
$ws.Range("Q940").Value = 15
$ws.Range("R940").Value = $clasif
